# promo_test_data.xlsx — test data refresh
# (1) Replace the stale admin e-mail used across the TC rows (shared string
#     "zwshashank.agrawal@teampureplay.com" -> "admin@sunrise.com") so every
#     cell that referenced it (B3:B69, the username_admin column) is updated
#     in one shot, same as re-pointing the shared-string table entry.
# (2) Move the saved cursor/selection on Sheet1 from D78 to C73.
#
# NOTE: the source diff also touched the Microsoft-only
# <x15ac:absPath .../> hint inside <mc:AlternateContent> (workbook.xml) and
# the <sheetView topLeftCell="..."> scroll anchor. Neither is reachable from
# the Excel object model (there is no VBA/COM property for the absPath
# extension, and Excel only persists topLeftCell when panes are frozen/split,
# which would change the workbook's semantics) so they are left as-is here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- (1) username_admin column: swap the old e-mail for the new one -------
# xlWhole (1) / xlByRows (1) / forward search / not case-sensitive, mirrors
# Excel's Ctrl+H "Replace All" across every cell that holds the old address
# (the whole username_admin column, B3:B69) in a single pass.
$ws.Cells.Replace(
    "zwshashank.agrawal@teampureplay.com",
    "admin@sunrise.com",
    1, 1, $false, $false, $false, $false) | Out-Null

# --- (2) Update the active selection on the sheet -------------------------
$ws.Activate() | Out-Null
$ws.Range("C73").Select() | Out-Null
